$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay text (matches source data formatting)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '72.995.01'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '3.957.81'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '612.76'
$ws.Range('E5').Value = '  +13.42%  '
$ws.Range('D6').Value = '165.48'
$ws.Range('E6').Value = '  +9.00%  '
$ws.Range('D7').Value = '0.679'
$ws.Range('E7').Value = '  -2.92%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').Value = '0.183'
$ws.Range('E10').Value = '  +6.51%  '
$ws.Range('D11').Value = '55.93'
$ws.Range('E11').Value = '  +3.65%  '
$ws.Range('D12').Value = '0.0000332'
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').Value = '11.11'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = '4.584.99'
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').Value = '3.956.98'
$ws.Range('E15').Value = '  -2.17%  '
$ws.Range('D16').Value = '1.25'
$ws.Range('E16').Value = '  +3.60%  '
$ws.Range('D17').Value = '14.06'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '20.45'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = '72.805.95'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = '438.43'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').Value = '4.88'
$ws.Range('E22').Value = '  +14.66%  '
$ws.Range('D23').Value = '95.75'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('D24').Value = '3.38'
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('D25').Value = '14.12'
$ws.Range('E25').Value = '  -3.74%  '
$ws.Range('D26').Value = '4.06'
$ws.Range('E26').Value = '  -5.47%  '
$ws.Range('D27').Value = '11.04'
$ws.Range('E27').Value = '  -2.25%  '
$ws.Range('D28').Value = '5.96'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '10.46'
$ws.Range('E29').Value = '  -2.98%  '
$ws.Range('D30').Value = '35.93'
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('D31').Value = '8.04'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('D32').Value = '13.61'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0000104'
$ws.Range('E33').Value = '  +18.02%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.129'
$ws.Range('E34').Value = '  -3.73%  '
$ws.Range('D35').Value = '47.53'
$ws.Range('E35').Value = '  -3.56%  '
$ws.Range('D36').Value = '70.27'
$ws.Range('E36').Value = '  +5.12%  '
$ws.Range('D37').Value = '640.34'
$ws.Range('E37').Value = '  -5.87%  '
$ws.Range('D38').Value = '0.431'
$ws.Range('E38').Value = '  -5.33%  '
$ws.Range('D39').Value = '3.47'
$ws.Range('E39').Value = '  +2.59%  '
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0482'
$ws.Range('E43').Value = '  -2.47%  '
$ws.Range('B44').Value = 'THORChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D44').Value = '10.67'
$ws.Range('E44').Value = '  -4.86%  '
$ws.Range('D45').Value = '3.19'
$ws.Range('E45').Value = '  -6.09%  '
$ws.Range('E46').Value = '  -1.66%  '
$ws.Range('E47').Value = '  +3.41%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '2.86'
$ws.Range('E49').Value = '  +24.62%  '
$ws.Range('D50').Value = '2.836.63'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('D51').Value = '149.83'
$ws.Range('E51').Value = '  +1.83%  '
